$wb = $excel.ActiveWorkbook

$sheetNames = @("Clinical Structures", "opt structures", "couch_structures")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $cols = $ws.Range("C1:F1")
    $cols.EntireColumn.Select() | Out-Null
    $cols.EntireColumn.Delete() | Out-Null
}
